$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 166669310
$ws.Cells.Item(40, 9).Value = 3974.5
$ws.Cells.Item(40, 10).Value = 500000000
$ws.Cells.Item(40, 11).Value = 3974.5
$ws.Cells.Item(40, 12).Value = 500000000
$ws.Cells.Item(40, 13).Value = -3799.5
$ws.Cells.Item(40, 14).Value = -500000350

# Row 43
$ws.Cells.Item(43, 8).Value = 2975.2727
$ws.Cells.Item(43, 9).Value = 3987.8333
$ws.Cells.Item(43, 10).Value = 1760.2
$ws.Cells.Item(43, 11).Value = 3987.8333
$ws.Cells.Item(43, 12).Value = 1760.2
$ws.Cells.Item(43, 13).Value = -3918.8333

# Row 86
$ws.Cells.Item(86, 8).Value = 7483.357
$ws.Cells.Item(86, 9).Value = 2977.6
$ws.Cells.Item(86, 10).Value = 9986.556
$ws.Cells.Item(86, 11).Value = 2977.6
$ws.Cells.Item(86, 12).Value = 9986.556
$ws.Cells.Item(86, 13).Value = -1854.6
$ws.Cells.Item(86, 14).Value = -12232.556

# Row 89
$ws.Cells.Item(89, 8).Value = 7483.357
$ws.Cells.Item(89, 9).Value = 2977.6
$ws.Cells.Item(89, 10).Value = 9986.556
$ws.Cells.Item(89, 11).Value = 14888
$ws.Cells.Item(89, 12).Value = 49932.78
$ws.Cells.Item(89, 13).Value = -9272
$ws.Cells.Item(89, 14).Value = -61164.78

# Row 96
$ws.Cells.Item(96, 8).Value = 1612911.2
$ws.Cells.Item(96, 9).Value = 2324.5
$ws.Cells.Item(96, 10).Value = 2901380.5
$ws.Cells.Item(96, 11).Value = 6973.5
$ws.Cells.Item(96, 12).Value = 8704141.5
$ws.Cells.Item(96, 13).Value = -5600.5
$ws.Cells.Item(96, 14).Value = -8706887.5

# Row 99
$ws.Cells.Item(99, 8).Value = 3429.2856
$ws.Cells.Item(99, 9).Value = 154.5
$ws.Cells.Item(99, 10).Value = 4739.2
$ws.Cells.Item(99, 11).Value = 463.5
$ws.Cells.Item(99, 12).Value = 14217.6
$ws.Cells.Item(99, 13).Value = 1034.5
$ws.Cells.Item(99, 14).Value = -17213.6

# Row 118
$ws.Cells.Item(118, 8).Value = 1428
$ws.Cells.Item(118, 9).Value = 1496.6666
$ws.Cells.Item(118, 10).Value = 1325
$ws.Cells.Item(118, 11).Value = 4489.9998
$ws.Cells.Item(118, 12).Value = 3975
$ws.Cells.Item(118, 13).Value = -2832.9998
$ws.Cells.Item(118, 14).Value = -7289

# Row 132
$ws.Cells.Item(132, 8).Value = 6142.84
$ws.Cells.Item(132, 9).Value = 3422.1177
$ws.Cells.Item(132, 10).Value = 11924.375
$ws.Cells.Item(132, 11).Value = 10266.3531
$ws.Cells.Item(132, 12).Value = 35773.125
$ws.Cells.Item(132, 13).Value = -7736.3531

# Row 137
$ws.Cells.Item(137, 8).Value = 2048
$ws.Cells.Item(137, 9).Value = 2070.0908
$ws.Cells.Item(137, 10).Value = 1999.4
$ws.Cells.Item(137, 11).Value = 6210.2724
$ws.Cells.Item(137, 12).Value = 5998.200000000001
$ws.Cells.Item(137, 13).Value = -3660.2724

$ws = $wb.Sheets.Item("ARM")
# Row 74
$ws.Cells.Item(74, 8).Value = 1988.6666
$ws.Cells.Item(74, 9).Value = 1380.6666
$ws.Cells.Item(74, 10).Value = 3508.6667
$ws.Cells.Item(74, 11).Value = 1380.6666
$ws.Cells.Item(74, 12).Value = 3508.6667
$ws.Cells.Item(74, 13).Value = -506.6666

# Row 77
$ws.Cells.Item(77, 8).Value = 1988.6666
$ws.Cells.Item(77, 9).Value = 1380.6666
$ws.Cells.Item(77, 10).Value = 3508.6667
$ws.Cells.Item(77, 11).Value = 6903.333000000001
$ws.Cells.Item(77, 12).Value = 17543.3335
$ws.Cells.Item(77, 13).Value = -2535.333000000001

# Row 102
$ws.Cells.Item(102, 8).Value = 1413.3334
$ws.Cells.Item(102, 9).Value = 1093.8
$ws.Cells.Item(102, 10).Value = 3011
$ws.Cells.Item(102, 11).Value = 1093.8
$ws.Cells.Item(102, 12).Value = 3011
$ws.Cells.Item(102, 13).Value = 528.2

# Row 122
$ws.Cells.Item(122, 8).Value = 3350.111
$ws.Cells.Item(122, 9).Value = 3350.111
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 10050.333
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -7600.332999999999
$ws.Cells.Item(122, 14).Value = ""

# Row 132
$ws.Cells.Item(132, 8).Value = 5004955
$ws.Cells.Item(132, 9).Value = 4406.467
$ws.Cells.Item(132, 10).Value = 20006600
$ws.Cells.Item(132, 11).Value = 13219.401
$ws.Cells.Item(132, 12).Value = 60019800
$ws.Cells.Item(132, 13).Value = -10689.401
$ws.Cells.Item(132, 14).Value = -60024860

$ws = $wb.Sheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 1894.6111
$ws.Cells.Item(94, 9).Value = 1996.1666
$ws.Cells.Item(94, 10).Value = 1691.5
$ws.Cells.Item(94, 11).Value = 1996.1666
$ws.Cells.Item(94, 12).Value = 1691.5
$ws.Cells.Item(94, 13).Value = -1545.1666
$ws.Cells.Item(94, 14).Value = -2593.5

# Row 105
$ws.Cells.Item(105, 8).Value = 717639.4399999999
$ws.Cells.Item(105, 9).Value = 1170155.8
$ws.Cells.Item(105, 10).Value = 6542.4287
$ws.Cells.Item(105, 11).Value = 1170155.8
$ws.Cells.Item(105, 12).Value = 6542.4287
$ws.Cells.Item(105, 13).Value = -1168408.8

# Row 124
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).Value = ""

# Row 125
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = ""

# Row 134
$ws.Cells.Item(134, 8).Value = 4002161
$ws.Cells.Item(134, 9).Value = 2175.1738
$ws.Cells.Item(134, 10).Value = 50001996
$ws.Cells.Item(134, 11).Value = 6525.5214
$ws.Cells.Item(134, 12).Value = 150005988
$ws.Cells.Item(134, 13).Value = -3990.5214

$ws = $wb.Sheets.Item("CRP")
# Row 48
$ws.Cells.Item(48, 8).Value = 89999
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 89999
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 89999
$ws.Cells.Item(48, 14).Value = -90951

# Row 99
$ws.Cells.Item(99, 8).Value = 19158.25
$ws.Cells.Item(99, 9).Value = 21609.428
$ws.Cells.Item(99, 10).Value = 2000
$ws.Cells.Item(99, 11).Value = 21609.428
$ws.Cells.Item(99, 12).Value = 2000
$ws.Cells.Item(99, 13).Value = -20111.428

# Row 105
$ws.Cells.Item(105, 8).Value = 1982.8572
$ws.Cells.Item(105, 9).Value = 1556.2142
$ws.Cells.Item(105, 10).Value = 2836.1428
$ws.Cells.Item(105, 11).Value = 1556.2142
$ws.Cells.Item(105, 12).Value = 2836.1428
$ws.Cells.Item(105, 13).Value = 190.7858000000001

# Row 107
$ws.Cells.Item(107, 8).Value = 1353.2693
$ws.Cells.Item(107, 9).Value = 948.4545000000001
$ws.Cells.Item(107, 10).Value = 3579.75
$ws.Cells.Item(107, 11).Value = 948.4545000000001
$ws.Cells.Item(107, 12).Value = 3579.75
$ws.Cells.Item(107, 13).Value = 971.5454999999999

# Row 125
$ws.Cells.Item(125, 8).Value = 92499.5
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 92499.5
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 92499.5
$ws.Cells.Item(125, 14).Value = -97419.5

# Row 126
$ws.Cells.Item(126, 8).Value = 19158.25
$ws.Cells.Item(126, 9).Value = 21609.428
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 64828.284
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -62358.284

$ws = $wb.Sheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 59.333332
$ws.Cells.Item(2, 9).Value = 80
$ws.Cells.Item(2, 10).Value = 55.2
$ws.Cells.Item(2, 11).Value = 480
$ws.Cells.Item(2, 12).Value = 331.2
$ws.Cells.Item(2, 13).Value = -367
$ws.Cells.Item(2, 14).Value = -557.2

# Row 11
$ws.Cells.Item(11, 8).Value = 5011.5557
$ws.Cells.Item(11, 9).Value = 1479.3334
$ws.Cells.Item(11, 10).Value = 12076
$ws.Cells.Item(11, 11).Value = 4438.0002
$ws.Cells.Item(11, 12).Value = 36228
$ws.Cells.Item(11, 13).Value = -4298.0002
$ws.Cells.Item(11, 14).Value = -36508

# Row 132
$ws.Cells.Item(132, 8).Value = 2904.5
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 2904.5
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 26140.5
$ws.Cells.Item(132, 13).Value = ""

$ws = $wb.Sheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 9304.632
$ws.Cells.Item(70, 9).Value = 3936
$ws.Cells.Item(70, 10).Value = 10311.25
$ws.Cells.Item(70, 11).Value = 3936
$ws.Cells.Item(70, 12).Value = 10311.25
$ws.Cells.Item(70, 13).Value = -3666

# Row 73
$ws.Cells.Item(73, 8).Value = 9304.632
$ws.Cells.Item(73, 9).Value = 3936
$ws.Cells.Item(73, 10).Value = 10311.25
$ws.Cells.Item(73, 11).Value = 3936
$ws.Cells.Item(73, 12).Value = 10311.25
$ws.Cells.Item(73, 13).Value = -3000

# Row 80
$ws.Cells.Item(80, 8).Value = 1839.6
$ws.Cells.Item(80, 9).Value = 1468.4
$ws.Cells.Item(80, 10).Value = 2582
$ws.Cells.Item(80, 11).Value = 1468.4
$ws.Cells.Item(80, 12).Value = 2582
$ws.Cells.Item(80, 13).Value = -470.4000000000001
$ws.Cells.Item(80, 14).Value = -4578

# Row 83
$ws.Cells.Item(83, 8).Value = 1839.6
$ws.Cells.Item(83, 9).Value = 1468.4
$ws.Cells.Item(83, 10).Value = 2582
$ws.Cells.Item(83, 11).Value = 7342
$ws.Cells.Item(83, 12).Value = 12910
$ws.Cells.Item(83, 13).Value = -2350
$ws.Cells.Item(83, 14).Value = -22894

# Row 102
$ws.Cells.Item(102, 8).Value = 3589.842
$ws.Cells.Item(102, 9).Value = 3592.8125
$ws.Cells.Item(102, 10).Value = 3574
$ws.Cells.Item(102, 11).Value = 3592.8125
$ws.Cells.Item(102, 12).Value = 3574
$ws.Cells.Item(102, 13).Value = -1970.8125

# Row 107
$ws.Cells.Item(107, 8).Value = 1249.1177
$ws.Cells.Item(107, 9).Value = 1217
$ws.Cells.Item(107, 10).Value = 1399
$ws.Cells.Item(107, 11).Value = 1217
$ws.Cells.Item(107, 12).Value = 1399
$ws.Cells.Item(107, 13).Value = 703
$ws.Cells.Item(107, 14).Value = -5239

# Row 113
$ws.Cells.Item(113, 8).Value = 1685900.2
$ws.Cells.Item(113, 9).Value = 2427.8572
$ws.Cells.Item(113, 10).Value = 4631977
$ws.Cells.Item(113, 11).Value = 2427.8572
$ws.Cells.Item(113, 12).Value = 4631977
$ws.Cells.Item(113, 13).Value = -257.8571999999999
$ws.Cells.Item(113, 14).Value = -4636317

# Row 126
$ws.Cells.Item(126, 8).Value = 3626.2222
$ws.Cells.Item(126, 9).Value = 3034.2144
$ws.Cells.Item(126, 10).Value = 5698.25
$ws.Cells.Item(126, 11).Value = 9102.643199999999
$ws.Cells.Item(126, 12).Value = 17094.75
$ws.Cells.Item(126, 13).Value = -6632.643199999999

$ws = $wb.Sheets.Item("LTW")
# Row 26
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = ""

# Row 32
$ws.Cells.Item(32, 8).Value = 2874.75
$ws.Cells.Item(32, 9).Value = 2874.75
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 2874.75
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -2557.75

# Row 40
$ws.Cells.Item(40, 8).Value = 4994.091
$ws.Cells.Item(40, 9).Value = 4493.5
$ws.Cells.Item(40, 10).Value = 10000
$ws.Cells.Item(40, 11).Value = 4493.5
$ws.Cells.Item(40, 12).Value = 10000
$ws.Cells.Item(40, 13).Value = -4357.5

# Row 43
$ws.Cells.Item(43, 8).Value = 8248.75
$ws.Cells.Item(43, 9).Value = 7666.6665
$ws.Cells.Item(43, 10).Value = 9995
$ws.Cells.Item(43, 11).Value = 7666.6665
$ws.Cells.Item(43, 12).Value = 9995
$ws.Cells.Item(43, 13).Value = -7473.6665
$ws.Cells.Item(43, 14).Value = -10381

# Row 61
$ws.Cells.Item(61, 8).Value = 83338970
$ws.Cells.Item(61, 9).Value = 166667620
$ws.Cells.Item(61, 10).Value = 10316.333
$ws.Cells.Item(61, 11).Value = 166667620
$ws.Cells.Item(61, 12).Value = 10316.333
$ws.Cells.Item(61, 13).Value = -166667418

# Row 68
$ws.Cells.Item(68, 8).Value = 20835334
$ws.Cells.Item(68, 9).Value = 41666664
$ws.Cells.Item(68, 10).Value = 4003
$ws.Cells.Item(68, 11).Value = 41666664
$ws.Cells.Item(68, 12).Value = 4003
$ws.Cells.Item(68, 13).Value = -41665915

# Row 71
$ws.Cells.Item(71, 8).Value = 20835334
$ws.Cells.Item(71, 9).Value = 41666664
$ws.Cells.Item(71, 10).Value = 4003
$ws.Cells.Item(71, 11).Value = 208333320
$ws.Cells.Item(71, 12).Value = 20015
$ws.Cells.Item(71, 13).Value = -208329576

# Row 82
$ws.Cells.Item(82, 8).Value = 4585.857
$ws.Cells.Item(82, 9).Value = 3038.25
$ws.Cells.Item(82, 10).Value = 6649.3335
$ws.Cells.Item(82, 11).Value = 3038.25
$ws.Cells.Item(82, 12).Value = 6649.3335
$ws.Cells.Item(82, 13).Value = -2677.25
$ws.Cells.Item(82, 14).Value = -7371.3335

# Row 85
$ws.Cells.Item(85, 8).Value = 4585.857
$ws.Cells.Item(85, 9).Value = 3038.25
$ws.Cells.Item(85, 10).Value = 6649.3335
$ws.Cells.Item(85, 11).Value = 3038.25
$ws.Cells.Item(85, 12).Value = 6649.3335
$ws.Cells.Item(85, 13).Value = -1790.25
$ws.Cells.Item(85, 14).Value = -9145.333500000001

# Row 93
$ws.Cells.Item(93, 8).Value = 2649318
$ws.Cells.Item(93, 9).Value = 2451.3333
$ws.Cells.Item(93, 10).Value = 18530518
$ws.Cells.Item(93, 11).Value = 2451.3333
$ws.Cells.Item(93, 12).Value = 18530518
$ws.Cells.Item(93, 13).Value = -1203.3333
$ws.Cells.Item(93, 14).Value = -18533014

# Row 110
$ws.Cells.Item(110, 8).Value = 71999.336
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 10).Value = 71999.336
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 12).Value = 71999.336
$ws.Cells.Item(110, 14).Value = -80179.336

# Row 113
$ws.Cells.Item(113, 8).Value = 83338970
$ws.Cells.Item(113, 9).Value = 166667620
$ws.Cells.Item(113, 10).Value = 10316.333
$ws.Cells.Item(113, 11).Value = 166667620
$ws.Cells.Item(113, 12).Value = 10316.333
$ws.Cells.Item(113, 13).Value = -166665450

# Row 122
$ws.Cells.Item(122, 8).Value = 3296.8572
$ws.Cells.Item(122, 9).Value = 3134.225
$ws.Cells.Item(122, 10).Value = 6549.5
$ws.Cells.Item(122, 11).Value = 9402.674999999999
$ws.Cells.Item(122, 12).Value = 19648.5
$ws.Cells.Item(122, 13).Value = -6952.674999999999
$ws.Cells.Item(122, 14).Value = -24548.5

$ws = $wb.Sheets.Item("WVR")
# Row 126
$ws.Cells.Item(126, 8).Value = 4124.6
$ws.Cells.Item(126, 9).Value = 4124.6
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 12373.8
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -9903.800000000001

# Row 132
$ws.Cells.Item(132, 8).Value = 242574.02
$ws.Cells.Item(132, 9).Value = 4215.297
$ws.Cells.Item(132, 10).Value = 1005321.9
$ws.Cells.Item(132, 11).Value = 12645.891
$ws.Cells.Item(132, 12).Value = 3015965.7
$ws.Cells.Item(132, 13).Value = -10115.891
